$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for column D cells being updated so that Excel does not
# auto-convert numeric-looking strings (e.g. "1.007") into floating point numbers
# and strip significant trailing zeros (e.g. "0.3620" -> 0.362).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.803.76'
$ws.Range("E2").Value = '  -2.61%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.776.65'
$ws.Range("E3").Value = '  -3.03%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.55%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.05'
$ws.Range("E5").Value = '  -1.17%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.32%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4279'
$ws.Range("E7").Value = '  -0.29%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3620'

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07186'
$ws.Range("E9").Value = '  -1.48%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8369'
$ws.Range("E10").Value = '  -3.65%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.26'
$ws.Range("E11").Value = '  -2.17%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.854.79'
$ws.Range("E12").Value = '  -2.67%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.257'
$ws.Range("E13").Value = '  -2.85%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.346'
$ws.Range("E14").Value = '  -2.93%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06790'
$ws.Range("E15").Value = '  -2.17%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.009'
$ws.Range("E16").Value = '  +0.55%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '79.47'
$ws.Range("E17").Value = '  -1.36%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008680'
$ws.Range("E18").Value = '  -2.81%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.005'
$ws.Range("E19").Value = '  +0.31%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.95'
$ws.Range("E20").Value = '  -3.37%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.913.61'
$ws.Range("E21").Value = '  -2.97%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.020'
$ws.Range("E22").Value = '  -2.62%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.02'
$ws.Range("E23").Value = '  +1.15%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.019.68'
$ws.Range("E24").Value = '  -3.79%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.917'
$ws.Range("E25").Value = '  -3.31%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.20'
$ws.Range("E26").Value = '  -0.96%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.19'
$ws.Range("E27").Value = '  -4.01%  '

# Row 28
$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.30'
$ws.Range("E28").Value = '  +0.81%  '

# Row 29
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.039'
$ws.Range("E29").Value = '  -2.25%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.625'
$ws.Range("E30").Value = '  -11.95%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08973'
$ws.Range("E31").Value = '  +1.28%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7221'
$ws.Range("E32").Value = '  -4.66%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.840'
$ws.Range("E33").Value = '  -4.67%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.327'
$ws.Range("E34").Value = '  -4.98%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.092'
$ws.Range("E35").Value = '  -4.22%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.004'
$ws.Range("E36").Value = '  +0.34%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.073'
$ws.Range("E37").Value = '  -1.42%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01895'
$ws.Range("E38").Value = '  -2.45%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05084'
$ws.Range("E39").Value = '  -4.76%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4913'
$ws.Range("E40").Value = '  -3.79%  '

# Row 41
$ws.Range("E41").Value = '  -3.78%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.546'
$ws.Range("E42").Value = '  -9.03%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.095'
$ws.Range("E43").Value = '  -7.57%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.920'
$ws.Range("E44").Value = '  -5.55%  '

# Row 45
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.74'
$ws.Range("E45").Value = '  -1.60%  '

# Row 46
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.004'
$ws.Range("E46").Value = '  +0.41%  '

# Row 47
$ws.Range("E47").Value = '  -4.31%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06223'
$ws.Range("E48").Value = '  -4.31%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4479'
$ws.Range("E49").Value = '  -4.51%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.572'
$ws.Range("E50").Value = '  -3.08%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.726'
$ws.Range("E51").Value = '  -0.82%  '
